$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as literal TEXT, even when it looks like
# a date (e.g. "2020-12-16"), without Excel's COM auto-detection turning it
# into a date serial number. We stage the text (quote-prefixed so Excel
# keeps it as text) in a scratch cell far outside the used range, copy just
# the value over with PasteSpecial(xlPasteValues), then wipe the scratch
# cell (format included) so it leaves no trace in the saved workbook.
function Set-TextValue($sheet, $row, $col, $text) {
    $scratch = $sheet.Cells.Item(1, 200)
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $sheet.Cells.Item($row, $col).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Sheet "cases_by_race": fix the running index on the 2020-12-14 block
# (rows 53-61 had reset to 0-8 instead of continuing 51-59) and append the
# new 2020-12-16 (as-of 2020-12-15) block as rows 62-70.
# ---------------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("cases_by_race")

$wsRace.Cells.Item(53, 1).Value = 51
$wsRace.Cells.Item(54, 1).Value = 52
$wsRace.Cells.Item(55, 1).Value = 53
$wsRace.Cells.Item(56, 1).Value = 54
$wsRace.Cells.Item(57, 1).Value = 55
$wsRace.Cells.Item(58, 1).Value = 56
$wsRace.Cells.Item(59, 1).Value = 57
$wsRace.Cells.Item(60, 1).Value = 58
$wsRace.Cells.Item(61, 1).Value = 59

# New rows 62-70 should carry the same "index column" styling (bold, border,
# centred) as the rest of column A; copy it across before writing values.
$wsRace.Range("A53:A61").Copy()
$wsRace.Range("A62:A70").PasteSpecial(-4122)

$raceRows = @(
    @{ Row = 62; A = 0; B = ""; E = 4 },
    @{ Row = 63; A = 1; B = "American Indian or Alaska Native"; E = 52 },
    @{ Row = 64; A = 2; B = "Asian"; E = 238 },
    @{ Row = 65; A = 3; B = "Black or African American"; E = 1397 },
    @{ Row = 66; A = 4; B = "Native Hawaiian or Other Pacific Islander"; E = 10 },
    @{ Row = 67; A = 5; B = "Not disclosed"; E = 1679 },
    @{ Row = 68; A = 6; B = "Other Race"; E = 379 },
    @{ Row = 69; A = 7; B = "Two or more"; E = 125 },
    @{ Row = 70; A = 8; B = "White"; E = 13427 }
)

foreach ($r in $raceRows) {
    $row = $r.Row
    $wsRace.Cells.Item($row, 1).Value = $r.A
    if ($r.B -ne "") {
        $wsRace.Cells.Item($row, 2).Value = $r.B
    }
    Set-TextValue $wsRace $row 3 "2020-12-16"
    Set-TextValue $wsRace $row 4 "2020-12-15"
    $wsRace.Cells.Item($row, 5).Value = $r.E
}

# ---------------------------------------------------------------------------
# Sheet "cases_by_ethnicity": same fix for the 2020-12-14 block (rows 23-25,
# should continue 21-23) and append the new 2020-12-16 block as rows 26-28.
# ---------------------------------------------------------------------------
$wsEth = $wb.Worksheets.Item("cases_by_ethnicity")

$wsEth.Cells.Item(23, 1).Value = 21
$wsEth.Cells.Item(24, 1).Value = 22
$wsEth.Cells.Item(25, 1).Value = 23

$wsEth.Range("A23:A25").Copy()
$wsEth.Range("A26:A28").PasteSpecial(-4122)

$ethRows = @(
    @{ Row = 26; A = 0; B = "Hispanic or Latino"; E = 370 },
    @{ Row = 27; A = 1; B = "Not Hispanic or Latino"; E = 13332 },
    @{ Row = 28; A = 2; B = "unknown"; E = 3609 }
)

foreach ($r in $ethRows) {
    $row = $r.Row
    $wsEth.Cells.Item($row, 1).Value = $r.A
    $wsEth.Cells.Item($row, 2).Value = $r.B
    Set-TextValue $wsEth $row 3 "2020-12-16"
    Set-TextValue $wsEth $row 4 "2020-12-15"
    $wsEth.Cells.Item($row, 5).Value = $r.E
}
